$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.948.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.387.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.87%  "
$ws.Range("E7").Value = "  +2.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.375.35"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  +10.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.632"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.72%  "
$ws.Range("E13").Value = "  +5.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.917.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("E17").Value = "  +3.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.379.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "64.792.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.991"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +14.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.13%  "
$ws.Range("E28").Value = "  +2.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "571.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "61.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.53%  "
$ws.Range("E35").Value = "  +2.82%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.92%  "
$ws.Range("E38").Value = "  -3.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0741"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.369"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.53%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.083.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.96%  "
$ws.Range("E45").Value = "  +5.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.135"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.07%  "
$ws.Range("E47").Value = "  +2.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.76%  "
$ws.Range("E51").Value = "  +4.66%  "
